{"js": "// Paragraph-level text replacements derived from the canonical OOXML diff.\n// Each entry is the full text of one <w:p> run (a unique, exact match in the\n// document) mapped to its replacement text.\nconst replacements = [\n  { oldText: \"Obvin\u011bn\u00fd [[PERSON_2]] ([[BIRTH_ID_1]], OP \u010d.[[PHONE_1]]) je st\u00edh\u00e1n pro trestn\u00fd \u010din loupe\u017ee dle \u00a7 173 tr. z\u00e1kon\u00edku. Ob\u017ealovan\u00fd [[PERSON_2]] se k \u010dinu doznal. Po\u0161kozen\u00fd [[PERSON_3]] ([[BIRTH_ID_2]], bytem Poln\u00ed 45, N\u00e1chod 547 01, [[PHONE_2]]) utrp\u011bl t\u011b\u017ekou \u00fajmu na zdrav\u00ed. Sv\u011bdkem byl Ing. [[PERSON_1]] ([[BIRTH_ID_3]], email: [[EMAIL_1]]), kter\u00fd potvrdil pr\u016fb\u011bh \u00fatoku. Pachatel [[PERSON_2]] byl odsouzen k trestu odn\u011bt\u00ed svobody 8 let.\", newText: \"Obvin\u011bn\u00fd [[PERSON_2]] ([[BIRTH_ID_1]], [[ID_CARD_1]]) je st\u00edh\u00e1n pro trestn\u00fd \u010din loupe\u017ee dle \u00a7 173 tr. z\u00e1kon\u00edku. Ob\u017ealovan\u00fd [[PERSON_2]] se k \u010dinu doznal. Po\u0161kozen\u00fd [[PERSON_3]] ([[BIRTH_ID_2]], [[ADDRESS_1]], [[PHONE_1]]) utrp\u011bl t\u011b\u017ekou \u00fajmu na zdrav\u00ed. Sv\u011bdkem byl Ing. [[PERSON_1]] ([[BIRTH_ID_3]], email: [[EMAIL_1]]), kter\u00fd potvrdil pr\u016fb\u011bh \u00fatoku. Pachatel [[PERSON_2]] byl odsouzen k trestu odn\u011bt\u00ed svobody 8 let.\" },\n  { oldText: \"Zadr\u017een\u00e1 [[PERSON_5]] ([[BIRTH_ID_4]], \u0158P \u010d. [[ID_CARD_1]], bytem Slezsk\u00e1 89, Opava 746 01) je vy\u0161et\u0159ov\u00e1na pro podvod. Kontrolovan\u00e1 [[PERSON_5]] m\u011bla ve spr\u00e1v\u011b [[BANK_1]] s podez\u0159el\u00fdmi transakcemi. St\u011b\u017eovatel Mgr. [[PERSON_4]] ([[BIRTH_ID_5]], [[ICO_1]]) podal trestn\u00ed ozn\u00e1men\u00ed. Propu\u0161t\u011bn\u00e1 [[PERSON_5]] byla po v\u00fdslechu propu\u0161t\u011bna na svobodu.\", newText: \"Zadr\u017een\u00e1 [[PERSON_5]] ([[BIRTH_ID_4]], \u0158P \u010d. [[ID_CARD_2]], [[ADDRESS_2]]) je vy\u0161et\u0159ov\u00e1na pro podvod. Kontrolovan\u00e1 [[PERSON_5]] m\u011bla ve spr\u00e1v\u011b [[BANK_1]] s podez\u0159el\u00fdmi transakcemi. St\u011b\u017eovatel Mgr. [[PERSON_4]] ([[BIRTH_ID_5]], [[ICO_1]]) podal trestn\u00ed ozn\u00e1men\u00ed. Propu\u0161t\u011bn\u00e1 [[PERSON_5]] byla po v\u00fdslechu propu\u0161t\u011bna na svobodu.\" },\n  { oldText: \"Dlu\u017en\u00edk [[PERSON_7]] ([[BIRTH_ID_6]], bytem Lesn\u00ed 67/3, Jihlava 586 01, [[PHONE_3]]) m\u00e1 evidov\u00e1no 9 exekuc\u00ed v celkov\u00e9 v\u00fd\u0161i 1 456 800 K\u010d. V\u011b\u0159itelem je spole\u010dnost CrediFin a.s., [[ICO_2]]. Dlu\u017enice [[PERSON_7]] po\u017e\u00e1dala o oddlu\u017een\u00ed. Ru\u010ditelem byl Bc. [[PERSON_6]] ([[BIRTH_ID_7]], [[BANK_2]]), kter\u00fd za ni zaplatil \u010d\u00e1stku 240 000 K\u010d.\", newText: \"Dlu\u017en\u00edk [[PERSON_7]] ([[BIRTH_ID_6]], [[ADDRESS_3]], [[PHONE_2]]) m\u00e1 evidov\u00e1no 9 exekuc\u00ed v celkov\u00e9 v\u00fd\u0161i 1 456 800 K\u010d. V\u011b\u0159itelem je spole\u010dnost CrediFin a.s., [[ICO_2]]. Dlu\u017enice [[PERSON_7]] po\u017e\u00e1dala o oddlu\u017een\u00ed. Ru\u010ditelem byl Bc. [[PERSON_6]] ([[BIRTH_ID_7]], [[BANK_2]]), kter\u00fd za ni zaplatil \u010d\u00e1stku 240 000 K\u010d.\" },\n  { oldText: \"\u017dadatel [[PERSON_8]] ([[BIRTH_ID_8]], I\u010cO podnikatele: [[ICO_3]], bytem Zahradn\u00ed 45, Hradec Kr\u00e1lov\u00e9 500 02, email: [[EMAIL_2]], [[PHONE_4]]) po\u017e\u00e1dal o \u00fav\u011br ve v\u00fd\u0161i 3 500 000 K\u010d. \u017dadateli [[PERSON_8]] byl \u00fav\u011br schv\u00e1len s \u00farokovou sazbou 6,25 % p.a. Klientem banky je ji\u017e od roku 2018. Poji\u0161t\u011bnec [[PERSON_8]] m\u00e1 \u017eivotn\u00ed poji\u0161t\u011bn\u00ed ve v\u00fd\u0161i 5 000 000 K\u010d u \u010cesk\u00e1 poji\u0161\u0165ovna a.s.\", newText: \"\u017dadatel [[PERSON_8]] ([[BIRTH_ID_8]], I\u010cO podnikatele: [[ICO_3]], [[ADDRESS_4]], email: [[EMAIL_2]], [[PHONE_3]]) po\u017e\u00e1dal o \u00fav\u011br ve v\u00fd\u0161i 3 500 000 K\u010d. \u017dadateli [[PERSON_8]] byl \u00fav\u011br schv\u00e1len s \u00farokovou sazbou 6,25 % p.a. Klientem banky je ji\u017e od roku 2018. Poji\u0161t\u011bnec [[PERSON_8]] m\u00e1 \u017eivotn\u00ed poji\u0161t\u011bn\u00ed ve v\u00fd\u0161i 5 000 000 K\u010d u \u010cesk\u00e1 poji\u0161\u0165ovna a.s.\" },\n  { oldText: \"Pacient prof. MUDr. [[PERSON_9]], CSc. ([[BIRTH_ID_9]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_5]], [[ADDRESS_1]] 602 00) byl hospitalizov\u00e1n pro infarkt myokardu I21.9. Nemocn\u00fd [[PERSON_9]] podstoupil urgentn\u00ed kardiochirurgick\u00fd z\u00e1krok. O\u0161et\u0159uj\u00edc\u00ed l\u00e9ka\u0159 MUDr. [[PERSON_10]], Ph.D. (registrace \u010cLK: 89012) indikoval n\u00e1slednou rehabilitaci. Vy\u0161et\u0159en\u00fd pacient [[PERSON_9]] byl po 14 dnech propu\u0161t\u011bn do dom\u00e1c\u00ed p\u00e9\u010de.\", newText: \"Pacient prof. MUDr. [[PERSON_9]], CSc. ([[BIRTH_ID_9]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_4]], [[ADDRESS_5]]) byl hospitalizov\u00e1n pro infarkt myokardu I21.9. Nemocn\u00fd [[PERSON_9]] podstoupil urgentn\u00ed kardiochirurgick\u00fd z\u00e1krok. O\u0161et\u0159uj\u00edc\u00ed l\u00e9ka\u0159 MUDr. [[PERSON_10]], Ph.D. (registrace \u010cLK: 89012) indikoval n\u00e1slednou rehabilitaci. Vy\u0161et\u0159en\u00fd pacient [[PERSON_9]] byl po 14 dnech propu\u0161t\u011bn do dom\u00e1c\u00ed p\u00e9\u010de.\" },\n  { oldText: \"Pacientka Ing. [[PERSON_11]] ([[BIRTH_ID_10]], poji\u0161t\u011bnka OZP[[PHONE_6]], bytem L\u00e1ze\u0148sk\u00e1 34, Teplice 415 01, [[PHONE_7]]) m\u00e1 diagnostikovanou roztrou\u0161enou skler\u00f3zu G35. L\u00e9\u010den\u00e1 [[PERSON_11]] dost\u00e1v\u00e1 biologickou l\u00e9\u010dbu interferonem beta. Operovan\u00e1 [[PERSON_11]] podstoupila v minulosti tak\u00e9 operaci kardiostimul\u00e1toru. Invalidn\u00ed d\u016fchod II. stupn\u011b \u010din\u00ed 12 800 K\u010d m\u011bs\u00ed\u010dn\u011b.\", newText: \"Pacientka Ing. [[PERSON_11]] ([[BIRTH_ID_10]], poji\u0161t\u011bnka OZP[[PHONE_5]], [[ADDRESS_6]], [[PHONE_6]]) m\u00e1 diagnostikovanou roztrou\u0161enou skler\u00f3zu G35. L\u00e9\u010den\u00e1 [[PERSON_11]] dost\u00e1v\u00e1 biologickou l\u00e9\u010dbu interferonem beta. Operovan\u00e1 [[PERSON_11]] podstoupila v minulosti tak\u00e9 operaci kardiostimul\u00e1toru. Invalidn\u00ed d\u016fchod II. stupn\u011b \u010din\u00ed 12 800 K\u010d m\u011bs\u00ed\u010dn\u011b.\" },\n  { oldText: \"Rozveden\u00fd [[PERSON_13]] ([[BIRTH_ID_11]], bytem Poln\u00ed 67, Pardubice 530 02, email: [[EMAIL_3]]) m\u00e1 povinnost platit v\u00fd\u017eivn\u00e9 na nezletil\u00e9 [[PERSON_14]] ([[BIRTH_ID_12]]) a [[PERSON_15]] ([[BIRTH_ID_13]]) ve v\u00fd\u0161i 9 500 K\u010d m\u011bs\u00ed\u010dn\u011b na ka\u017ed\u00e9 d\u00edt\u011b. Pl\u00e1tce [[PERSON_13]] je v prodlen\u00ed s platbami celkem 57 000 K\u010d. P\u0159\u00edjemkyn\u00ed v\u00fd\u017eivn\u00e9ho je b\u00fdval\u00e1 man\u017eelka PhDr. [[PERSON_12]] ([[BIRTH_ID_14]], [[BANK_3]]).\", newText: \"Rozveden\u00fd [[PERSON_15]] ([[BIRTH_ID_11]], [[ADDRESS_7]], email: [[EMAIL_3]]) m\u00e1 povinnost platit v\u00fd\u017eivn\u00e9 na nezletil\u00e9 [[PERSON_14]] ([[BIRTH_ID_12]]) a [[PERSON_13]] ([[BIRTH_ID_13]]) ve v\u00fd\u0161i 9 500 K\u010d m\u011bs\u00ed\u010dn\u011b na ka\u017ed\u00e9 d\u00edt\u011b. Pl\u00e1tce [[PERSON_15]] je v prodlen\u00ed s platbami celkem 57 000 K\u010d. P\u0159\u00edjemkyn\u00ed v\u00fd\u017eivn\u00e9ho je b\u00fdval\u00e1 man\u017eelka PhDr. [[PERSON_12]] ([[BIRTH_ID_14]], [[BANK_3]]).\" },\n  { oldText: \"Opatrovn\u00edk JUDr. [[PERSON_16]] ([[BIRTH_ID_15]], registrace \u010cAK: 67890, bytem Kr\u00e1tk\u00e1 45, Olomouc 779 00) byl ustanoven pro [[PERSON_17]] ([[BIRTH_ID_16]], bytem Zahradn\u00ed 12, P\u0159erov 750 02), kter\u00e1 byla zbavena zp\u016fsobilosti k pr\u00e1vn\u00edm \u00fakon\u016fm. Opatrovance [[PERSON_17]] n\u00e1le\u017e\u00ed d\u016fchod 19 400 K\u010d m\u011bs\u00ed\u010dn\u011b a nemovitost v hodnot\u011b 6 200 000 K\u010d na LV \u010d. 7823.\", newText: \"Opatrovn\u00edk JUDr. [[PERSON_16]] ([[BIRTH_ID_15]], registrace \u010cAK: 67890, [[ADDRESS_8]]) byl ustanoven pro [[PERSON_17]] ([[BIRTH_ID_16]], [[ADDRESS_9]]), kter\u00e1 byla zbavena zp\u016fsobilosti k pr\u00e1vn\u00edm \u00fakon\u016fm. Opatrovance [[PERSON_17]] n\u00e1le\u017e\u00ed d\u016fchod 19 400 K\u010d m\u011bs\u00ed\u010dn\u011b a nemovitost v hodnot\u011b 6 200 000 K\u010d na LV \u010d. 7823.\" },\n  { oldText: \"Zam\u011bstnanec Mgr. [[PERSON_18]] ([[BIRTH_ID_17]], osobn\u00ed \u010d\u00edslo: 2023-HR-0789, [[ADDRESS_2]] 708 00, email: [[EMAIL_4]], [[PHONE_8]]) podal v\u00fdpov\u011b\u010f. Zam\u011bstnanci [[PERSON_18]] bylo vyplaceno odstupn\u00e9 ve v\u00fd\u0161i 186 400 K\u010d na [[BANK_4]]. Jeho n\u00e1stupcem je uchaze\u010d Bc. [[PERSON_19]] ([[BIRTH_ID_18]], bytem Kr\u00e1tk\u00e1 23, Hav\u00ed\u0159ov 736 01).\", newText: \"Zam\u011bstnanec Mgr. [[PERSON_18]] ([[BIRTH_ID_17]], osobn\u00ed \u010d\u00edslo: 2023-HR-0789, [[ADDRESS_10]], email: [[EMAIL_4]], [[PHONE_7]]) podal v\u00fdpov\u011b\u010f. Zam\u011bstnanci [[PERSON_18]] bylo vyplaceno odstupn\u00e9 ve v\u00fd\u0161i 186 400 K\u010d na [[BANK_4]]. Jeho n\u00e1stupcem je uchaze\u010d Bc. [[PERSON_19]] ([[BIRTH_ID_18]], [[ADDRESS_11]]).\" },\n  { oldText: \"D\u016fchodkyn\u011b [[PERSON_21]] ([[BIRTH_ID_19]], bytem Poln\u00ed 67, Fr\u00fddek-M\u00edstek 738 01, [[PHONE_9]]) pob\u00edr\u00e1 starobn\u00ed d\u016fchod ve v\u00fd\u0161i 17 800 K\u010d m\u011bs\u00ed\u010dn\u011b. Opr\u00e1vn\u011bn\u00e1 [[PERSON_20]] m\u00e1 n\u00e1rok tak\u00e9 na p\u0159\u00edsp\u011bvek na bydlen\u00ed 4 500 K\u010d. Dr\u017eitelka pr\u016fkazu ZTP m\u00e1 n\u00e1rok na slevy v MHD. Beneficientkou jej\u00edho \u017eivotn\u00edho poji\u0161t\u011bn\u00ed je dcera Ing. [[PERSON_20]] ([[BIRTH_ID_20]]).\", newText: \"D\u016fchodkyn\u011b [[PERSON_21]] ([[BIRTH_ID_19]], [[ADDRESS_12]], Fr\u00fddek-M\u00edstek 738 01, [[PHONE_8]]) pob\u00edr\u00e1 starobn\u00ed d\u016fchod ve v\u00fd\u0161i 17 800 K\u010d m\u011bs\u00ed\u010dn\u011b. Opr\u00e1vn\u011bn\u00e1 [[PERSON_20]] m\u00e1 n\u00e1rok tak\u00e9 na p\u0159\u00edsp\u011bvek na bydlen\u00ed 4 500 K\u010d. Dr\u017eitelka pr\u016fkazu ZTP m\u00e1 n\u00e1rok na slevy v MHD. Beneficientkou jej\u00edho \u017eivotn\u00edho poji\u0161t\u011bn\u00ed je dcera Ing. [[PERSON_20]] ([[BIRTH_ID_20]]).\" },\n  { oldText: \"\u017d\u00e1kyn\u011b [[PERSON_23]] ([[BIRTH_ID_21]], bytem L\u00e1ze\u0148sk\u00e1 45, Karlovy Vary 360 01) byla p\u0159ijata na Gymn\u00e1zium Karlovy Vary. Studentka [[PERSON_22]] z\u00edskala soci\u00e1ln\u00ed stipendium 3 500 K\u010d m\u011bs\u00ed\u010dn\u011b. Z\u00e1konnou z\u00e1stupkyn\u00ed je Mgr. [[PERSON_22]] ([[BIRTH_ID_14]], email: [[EMAIL_5]], [[PHONE_10]]). P\u0159ijat\u00e1 \u017e\u00e1kyn\u011b [[PERSON_23]] za\u010d\u00edn\u00e1 doch\u00e1zku od 1. 9. 2026.\", newText: \"\u017d\u00e1kyn\u011b [[PERSON_23]] ([[BIRTH_ID_21]], [[ADDRESS_13]]) byla p\u0159ijata na Gymn\u00e1zium Karlovy Vary. Studentka [[PERSON_22]] z\u00edskala soci\u00e1ln\u00ed stipendium 3 500 K\u010d m\u011bs\u00ed\u010dn\u011b. Z\u00e1konnou z\u00e1stupkyn\u00ed je Mgr. [[PERSON_22]] ([[BIRTH_ID_14]], email: [[EMAIL_5]], [[PHONE_9]]). P\u0159ijat\u00e1 \u017e\u00e1kyn\u011b [[PERSON_23]] za\u010d\u00edn\u00e1 doch\u00e1zku od 1. 9. 2026.\" },\n  { oldText: \"\u2022 Pacientka Ing. [[PERSON_26]] ([[BIRTH_ID_24]], [[PHONE_3]]) - l\u00e9\u010dba rakoviny vaje\u010dn\u00edk\u016f C56.9\", newText: \"\u2022 Pacientka Ing. [[PERSON_26]] ([[BIRTH_ID_24]], [[PHONE_2]]) - l\u00e9\u010dba rakoviny vaje\u010dn\u00edk\u016f C56.9\" },\n  { oldText: \"\u2022 Poji\u0161t\u011bnka PhDr. [[PERSON_31]] ([[BIRTH_ID_28]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_11]]) - hospitalizace pro pneumonii\", newText: \"\u2022 Poji\u0161t\u011bnka PhDr. [[PERSON_31]] ([[BIRTH_ID_28]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_10]]) - hospitalizace pro pneumonii\" },\n  { oldText: \"\u2022 Studentka [[PERSON_37]] ([[BIRTH_ID_34]], [[ADDRESS_3]] 602 00) - UK Praha, 2. ro\u010dn\u00edk\", newText: \"\u2022 Studentka [[PERSON_37]] ([[BIRTH_ID_34]], [[ADDRESS_14]]) - UK Praha, 2. ro\u010dn\u00edk\" },\n  { oldText: \"\u2022 Sv\u011bdkyn\u011b PharmDr. [[PERSON_39]] ([[BIRTH_ID_36]], [[PHONE_9]]) - v\u00fdpov\u011b\u010f v trestn\u00ed v\u011bci\", newText: \"\u2022 Sv\u011bdkyn\u011b PharmDr. [[PERSON_39]] ([[BIRTH_ID_36]], [[PHONE_8]]) - v\u00fdpov\u011b\u010f v trestn\u00ed v\u011bci\" },\n];\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load('items');\nawait context.sync();\n\nfor (let i = 0; i < paras.items.length; i++) {\n  paras.items[i].load('text');\n}\nawait context.sync();\n\nconst remaining = new Set(replacements.map((_, idx) => idx));\nlet appliedCount = 0;\nfor (let i = 0; i < paras.items.length; i++) {\n  const para = paras.items[i];\n  const text = para.text;\n  for (let ri = 0; ri < replacements.length; ri++) {\n    if (!remaining.has(ri)) continue;\n    const rep = replacements[ri];\n    if (text === rep.oldText) {\n      para.insertText(rep.newText, Word.InsertLocation.replace);\n      remaining.delete(ri);\n      appliedCount++;\n      break;\n    }\n  }\n}\nawait context.sync();\n\nif (remaining.size > 0) {\n  const missing = Array.from(remaining).map((ri) => replacements[ri].oldText.substring(0, 60));\n  throw new Error('Failed to locate ' + remaining.size + ' paragraph(s) for replacement: ' + JSON.stringify(missing));\n}\n\nreturn 'applied=' + appliedCount;", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = 'Obvin\u011bn\u00fd [[PERSON_2]] ([[BIRTH_ID_1]], OP \u010d.[[PHONE_1]]) je st\u00edh\u00e1n pro trestn\u00fd \u010din loupe\u017ee dle \u00a7 173 tr. z\u00e1kon\u00edku. Ob\u017ealovan\u00fd [[PERSON_2]] se k \u010dinu doznal. Po\u0161kozen\u00fd [[PERSON_3]] ([[BIRTH_ID_2]], bytem Poln\u00ed 45, N\u00e1chod 547 01, [[PHONE_2]]) utrp\u011bl t\u011b\u017ekou \u00fajmu na zdrav\u00ed. Sv\u011bdkem byl Ing. [[PERSON_1]] ([[BIRTH_ID_3]], email: [[EMAIL_1]]), kter\u00fd potvrdil pr\u016fb\u011bh \u00fatoku. Pachatel [[PERSON_2]] byl odsouzen k trestu odn\u011bt\u00ed svobody 8 let.'; New = 'Obvin\u011bn\u00fd [[PERSON_2]] ([[BIRTH_ID_1]], [[ID_CARD_1]]) je st\u00edh\u00e1n pro trestn\u00fd \u010din loupe\u017ee dle \u00a7 173 tr. z\u00e1kon\u00edku. Ob\u017ealovan\u00fd [[PERSON_2]] se k \u010dinu doznal. Po\u0161kozen\u00fd [[PERSON_3]] ([[BIRTH_ID_2]], [[ADDRESS_1]], [[PHONE_1]]) utrp\u011bl t\u011b\u017ekou \u00fajmu na zdrav\u00ed. Sv\u011bdkem byl Ing. [[PERSON_1]] ([[BIRTH_ID_3]], email: [[EMAIL_1]]), kter\u00fd potvrdil pr\u016fb\u011bh \u00fatoku. Pachatel [[PERSON_2]] byl odsouzen k trestu odn\u011bt\u00ed svobody 8 let.' },\n    @{ Old = 'Zadr\u017een\u00e1 [[PERSON_5]] ([[BIRTH_ID_4]], \u0158P \u010d. [[ID_CARD_1]], bytem Slezsk\u00e1 89, Opava 746 01) je vy\u0161et\u0159ov\u00e1na pro podvod. Kontrolovan\u00e1 [[PERSON_5]] m\u011bla ve spr\u00e1v\u011b [[BANK_1]] s podez\u0159el\u00fdmi transakcemi. St\u011b\u017eovatel Mgr. [[PERSON_4]] ([[BIRTH_ID_5]], [[ICO_1]]) podal trestn\u00ed ozn\u00e1men\u00ed. Propu\u0161t\u011bn\u00e1 [[PERSON_5]] byla po v\u00fdslechu propu\u0161t\u011bna na svobodu.'; New = 'Zadr\u017een\u00e1 [[PERSON_5]] ([[BIRTH_ID_4]], \u0158P \u010d. [[ID_CARD_2]], [[ADDRESS_2]]) je vy\u0161et\u0159ov\u00e1na pro podvod. Kontrolovan\u00e1 [[PERSON_5]] m\u011bla ve spr\u00e1v\u011b [[BANK_1]] s podez\u0159el\u00fdmi transakcemi. St\u011b\u017eovatel Mgr. [[PERSON_4]] ([[BIRTH_ID_5]], [[ICO_1]]) podal trestn\u00ed ozn\u00e1men\u00ed. Propu\u0161t\u011bn\u00e1 [[PERSON_5]] byla po v\u00fdslechu propu\u0161t\u011bna na svobodu.' },\n    @{ Old = 'Dlu\u017en\u00edk [[PERSON_7]] ([[BIRTH_ID_6]], bytem Lesn\u00ed 67/3, Jihlava 586 01, [[PHONE_3]]) m\u00e1 evidov\u00e1no 9 exekuc\u00ed v celkov\u00e9 v\u00fd\u0161i 1 456 800 K\u010d. V\u011b\u0159itelem je spole\u010dnost CrediFin a.s., [[ICO_2]]. Dlu\u017enice [[PERSON_7]] po\u017e\u00e1dala o oddlu\u017een\u00ed. Ru\u010ditelem byl Bc. [[PERSON_6]] ([[BIRTH_ID_7]], [[BANK_2]]), kter\u00fd za ni zaplatil \u010d\u00e1stku 240 000 K\u010d.'; New = 'Dlu\u017en\u00edk [[PERSON_7]] ([[BIRTH_ID_6]], [[ADDRESS_3]], [[PHONE_2]]) m\u00e1 evidov\u00e1no 9 exekuc\u00ed v celkov\u00e9 v\u00fd\u0161i 1 456 800 K\u010d. V\u011b\u0159itelem je spole\u010dnost CrediFin a.s., [[ICO_2]]. Dlu\u017enice [[PERSON_7]] po\u017e\u00e1dala o oddlu\u017een\u00ed. Ru\u010ditelem byl Bc. [[PERSON_6]] ([[BIRTH_ID_7]], [[BANK_2]]), kter\u00fd za ni zaplatil \u010d\u00e1stku 240 000 K\u010d.' },\n    @{ Old = '\u017dadatel [[PERSON_8]] ([[BIRTH_ID_8]], I\u010cO podnikatele: [[ICO_3]], bytem Zahradn\u00ed 45, Hradec Kr\u00e1lov\u00e9 500 02, email: [[EMAIL_2]], [[PHONE_4]]) po\u017e\u00e1dal o \u00fav\u011br ve v\u00fd\u0161i 3 500 000 K\u010d. \u017dadateli [[PERSON_8]] byl \u00fav\u011br schv\u00e1len s \u00farokovou sazbou 6,25 % p.a. Klientem banky je ji\u017e od roku 2018. Poji\u0161t\u011bnec [[PERSON_8]] m\u00e1 \u017eivotn\u00ed poji\u0161t\u011bn\u00ed ve v\u00fd\u0161i 5 000 000 K\u010d u \u010cesk\u00e1 poji\u0161\u0165ovna a.s.'; New = '\u017dadatel [[PERSON_8]] ([[BIRTH_ID_8]], I\u010cO podnikatele: [[ICO_3]], [[ADDRESS_4]], email: [[EMAIL_2]], [[PHONE_3]]) po\u017e\u00e1dal o \u00fav\u011br ve v\u00fd\u0161i 3 500 000 K\u010d. \u017dadateli [[PERSON_8]] byl \u00fav\u011br schv\u00e1len s \u00farokovou sazbou 6,25 % p.a. Klientem banky je ji\u017e od roku 2018. Poji\u0161t\u011bnec [[PERSON_8]] m\u00e1 \u017eivotn\u00ed poji\u0161t\u011bn\u00ed ve v\u00fd\u0161i 5 000 000 K\u010d u \u010cesk\u00e1 poji\u0161\u0165ovna a.s.' },\n    @{ Old = 'Pacient prof. MUDr. [[PERSON_9]], CSc. ([[BIRTH_ID_9]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_5]], [[ADDRESS_1]] 602 00) byl hospitalizov\u00e1n pro infarkt myokardu I21.9. Nemocn\u00fd [[PERSON_9]] podstoupil urgentn\u00ed kardiochirurgick\u00fd z\u00e1krok. O\u0161et\u0159uj\u00edc\u00ed l\u00e9ka\u0159 MUDr. [[PERSON_10]], Ph.D. (registrace \u010cLK: 89012) indikoval n\u00e1slednou rehabilitaci. Vy\u0161et\u0159en\u00fd pacient [[PERSON_9]] byl po 14 dnech propu\u0161t\u011bn do dom\u00e1c\u00ed p\u00e9\u010de.'; New = 'Pacient prof. MUDr. [[PERSON_9]], CSc. ([[BIRTH_ID_9]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_4]], [[ADDRESS_5]]) byl hospitalizov\u00e1n pro infarkt myokardu I21.9. Nemocn\u00fd [[PERSON_9]] podstoupil urgentn\u00ed kardiochirurgick\u00fd z\u00e1krok. O\u0161et\u0159uj\u00edc\u00ed l\u00e9ka\u0159 MUDr. [[PERSON_10]], Ph.D. (registrace \u010cLK: 89012) indikoval n\u00e1slednou rehabilitaci. Vy\u0161et\u0159en\u00fd pacient [[PERSON_9]] byl po 14 dnech propu\u0161t\u011bn do dom\u00e1c\u00ed p\u00e9\u010de.' },\n    @{ Old = 'Pacientka Ing. [[PERSON_11]] ([[BIRTH_ID_10]], poji\u0161t\u011bnka OZP[[PHONE_6]], bytem L\u00e1ze\u0148sk\u00e1 34, Teplice 415 01, [[PHONE_7]]) m\u00e1 diagnostikovanou roztrou\u0161enou skler\u00f3zu G35. L\u00e9\u010den\u00e1 [[PERSON_11]] dost\u00e1v\u00e1 biologickou l\u00e9\u010dbu interferonem beta. Operovan\u00e1 [[PERSON_11]] podstoupila v minulosti tak\u00e9 operaci kardiostimul\u00e1toru. Invalidn\u00ed d\u016fchod II. stupn\u011b \u010din\u00ed 12 800 K\u010d m\u011bs\u00ed\u010dn\u011b.'; New = 'Pacientka Ing. [[PERSON_11]] ([[BIRTH_ID_10]], poji\u0161t\u011bnka OZP[[PHONE_5]], [[ADDRESS_6]], [[PHONE_6]]) m\u00e1 diagnostikovanou roztrou\u0161enou skler\u00f3zu G35. L\u00e9\u010den\u00e1 [[PERSON_11]] dost\u00e1v\u00e1 biologickou l\u00e9\u010dbu interferonem beta. Operovan\u00e1 [[PERSON_11]] podstoupila v minulosti tak\u00e9 operaci kardiostimul\u00e1toru. Invalidn\u00ed d\u016fchod II. stupn\u011b \u010din\u00ed 12 800 K\u010d m\u011bs\u00ed\u010dn\u011b.' },\n    @{ Old = 'Rozveden\u00fd [[PERSON_13]] ([[BIRTH_ID_11]], bytem Poln\u00ed 67, Pardubice 530 02, email: [[EMAIL_3]]) m\u00e1 povinnost platit v\u00fd\u017eivn\u00e9 na nezletil\u00e9 [[PERSON_14]] ([[BIRTH_ID_12]]) a [[PERSON_15]] ([[BIRTH_ID_13]]) ve v\u00fd\u0161i 9 500 K\u010d m\u011bs\u00ed\u010dn\u011b na ka\u017ed\u00e9 d\u00edt\u011b. Pl\u00e1tce [[PERSON_13]] je v prodlen\u00ed s platbami celkem 57 000 K\u010d. P\u0159\u00edjemkyn\u00ed v\u00fd\u017eivn\u00e9ho je b\u00fdval\u00e1 man\u017eelka PhDr. [[PERSON_12]] ([[BIRTH_ID_14]], [[BANK_3]]).'; New = 'Rozveden\u00fd [[PERSON_15]] ([[BIRTH_ID_11]], [[ADDRESS_7]], email: [[EMAIL_3]]) m\u00e1 povinnost platit v\u00fd\u017eivn\u00e9 na nezletil\u00e9 [[PERSON_14]] ([[BIRTH_ID_12]]) a [[PERSON_13]] ([[BIRTH_ID_13]]) ve v\u00fd\u0161i 9 500 K\u010d m\u011bs\u00ed\u010dn\u011b na ka\u017ed\u00e9 d\u00edt\u011b. Pl\u00e1tce [[PERSON_15]] je v prodlen\u00ed s platbami celkem 57 000 K\u010d. P\u0159\u00edjemkyn\u00ed v\u00fd\u017eivn\u00e9ho je b\u00fdval\u00e1 man\u017eelka PhDr. [[PERSON_12]] ([[BIRTH_ID_14]], [[BANK_3]]).' },\n    @{ Old = 'Opatrovn\u00edk JUDr. [[PERSON_16]] ([[BIRTH_ID_15]], registrace \u010cAK: 67890, bytem Kr\u00e1tk\u00e1 45, Olomouc 779 00) byl ustanoven pro [[PERSON_17]] ([[BIRTH_ID_16]], bytem Zahradn\u00ed 12, P\u0159erov 750 02), kter\u00e1 byla zbavena zp\u016fsobilosti k pr\u00e1vn\u00edm \u00fakon\u016fm. Opatrovance [[PERSON_17]] n\u00e1le\u017e\u00ed d\u016fchod 19 400 K\u010d m\u011bs\u00ed\u010dn\u011b a nemovitost v hodnot\u011b 6 200 000 K\u010d na LV \u010d. 7823.'; New = 'Opatrovn\u00edk JUDr. [[PERSON_16]] ([[BIRTH_ID_15]], registrace \u010cAK: 67890, [[ADDRESS_8]]) byl ustanoven pro [[PERSON_17]] ([[BIRTH_ID_16]], [[ADDRESS_9]]), kter\u00e1 byla zbavena zp\u016fsobilosti k pr\u00e1vn\u00edm \u00fakon\u016fm. Opatrovance [[PERSON_17]] n\u00e1le\u017e\u00ed d\u016fchod 19 400 K\u010d m\u011bs\u00ed\u010dn\u011b a nemovitost v hodnot\u011b 6 200 000 K\u010d na LV \u010d. 7823.' },\n    @{ Old = 'Zam\u011bstnanec Mgr. [[PERSON_18]] ([[BIRTH_ID_17]], osobn\u00ed \u010d\u00edslo: 2023-HR-0789, [[ADDRESS_2]] 708 00, email: [[EMAIL_4]], [[PHONE_8]]) podal v\u00fdpov\u011b\u010f. Zam\u011bstnanci [[PERSON_18]] bylo vyplaceno odstupn\u00e9 ve v\u00fd\u0161i 186 400 K\u010d na [[BANK_4]]. Jeho n\u00e1stupcem je uchaze\u010d Bc. [[PERSON_19]] ([[BIRTH_ID_18]], bytem Kr\u00e1tk\u00e1 23, Hav\u00ed\u0159ov 736 01).'; New = 'Zam\u011bstnanec Mgr. [[PERSON_18]] ([[BIRTH_ID_17]], osobn\u00ed \u010d\u00edslo: 2023-HR-0789, [[ADDRESS_10]], email: [[EMAIL_4]], [[PHONE_7]]) podal v\u00fdpov\u011b\u010f. Zam\u011bstnanci [[PERSON_18]] bylo vyplaceno odstupn\u00e9 ve v\u00fd\u0161i 186 400 K\u010d na [[BANK_4]]. Jeho n\u00e1stupcem je uchaze\u010d Bc. [[PERSON_19]] ([[BIRTH_ID_18]], [[ADDRESS_11]]).' },\n    @{ Old = 'D\u016fchodkyn\u011b [[PERSON_21]] ([[BIRTH_ID_19]], bytem Poln\u00ed 67, Fr\u00fddek-M\u00edstek 738 01, [[PHONE_9]]) pob\u00edr\u00e1 starobn\u00ed d\u016fchod ve v\u00fd\u0161i 17 800 K\u010d m\u011bs\u00ed\u010dn\u011b. Opr\u00e1vn\u011bn\u00e1 [[PERSON_20]] m\u00e1 n\u00e1rok tak\u00e9 na p\u0159\u00edsp\u011bvek na bydlen\u00ed 4 500 K\u010d. Dr\u017eitelka pr\u016fkazu ZTP m\u00e1 n\u00e1rok na slevy v MHD. Beneficientkou jej\u00edho \u017eivotn\u00edho poji\u0161t\u011bn\u00ed je dcera Ing. [[PERSON_20]] ([[BIRTH_ID_20]]).'; New = 'D\u016fchodkyn\u011b [[PERSON_21]] ([[BIRTH_ID_19]], [[ADDRESS_12]], Fr\u00fddek-M\u00edstek 738 01, [[PHONE_8]]) pob\u00edr\u00e1 starobn\u00ed d\u016fchod ve v\u00fd\u0161i 17 800 K\u010d m\u011bs\u00ed\u010dn\u011b. Opr\u00e1vn\u011bn\u00e1 [[PERSON_20]] m\u00e1 n\u00e1rok tak\u00e9 na p\u0159\u00edsp\u011bvek na bydlen\u00ed 4 500 K\u010d. Dr\u017eitelka pr\u016fkazu ZTP m\u00e1 n\u00e1rok na slevy v MHD. Beneficientkou jej\u00edho \u017eivotn\u00edho poji\u0161t\u011bn\u00ed je dcera Ing. [[PERSON_20]] ([[BIRTH_ID_20]]).' },\n    @{ Old = '\u017d\u00e1kyn\u011b [[PERSON_23]] ([[BIRTH_ID_21]], bytem L\u00e1ze\u0148sk\u00e1 45, Karlovy Vary 360 01) byla p\u0159ijata na Gymn\u00e1zium Karlovy Vary. Studentka [[PERSON_22]] z\u00edskala soci\u00e1ln\u00ed stipendium 3 500 K\u010d m\u011bs\u00ed\u010dn\u011b. Z\u00e1konnou z\u00e1stupkyn\u00ed je Mgr. [[PERSON_22]] ([[BIRTH_ID_14]], email: [[EMAIL_5]], [[PHONE_10]]). P\u0159ijat\u00e1 \u017e\u00e1kyn\u011b [[PERSON_23]] za\u010d\u00edn\u00e1 doch\u00e1zku od 1. 9. 2026.'; New = '\u017d\u00e1kyn\u011b [[PERSON_23]] ([[BIRTH_ID_21]], [[ADDRESS_13]]) byla p\u0159ijata na Gymn\u00e1zium Karlovy Vary. Studentka [[PERSON_22]] z\u00edskala soci\u00e1ln\u00ed stipendium 3 500 K\u010d m\u011bs\u00ed\u010dn\u011b. Z\u00e1konnou z\u00e1stupkyn\u00ed je Mgr. [[PERSON_22]] ([[BIRTH_ID_14]], email: [[EMAIL_5]], [[PHONE_9]]). P\u0159ijat\u00e1 \u017e\u00e1kyn\u011b [[PERSON_23]] za\u010d\u00edn\u00e1 doch\u00e1zku od 1. 9. 2026.' },\n    @{ Old = '\u2022 Pacientka Ing. [[PERSON_26]] ([[BIRTH_ID_24]], [[PHONE_3]]) - l\u00e9\u010dba rakoviny vaje\u010dn\u00edk\u016f C56.9'; New = '\u2022 Pacientka Ing. [[PERSON_26]] ([[BIRTH_ID_24]], [[PHONE_2]]) - l\u00e9\u010dba rakoviny vaje\u010dn\u00edk\u016f C56.9' },\n    @{ Old = '\u2022 Poji\u0161t\u011bnka PhDr. [[PERSON_31]] ([[BIRTH_ID_28]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_11]]) - hospitalizace pro pneumonii'; New = '\u2022 Poji\u0161t\u011bnka PhDr. [[PERSON_31]] ([[BIRTH_ID_28]], \u010d\u00edslo poji\u0161t\u011bnce VZP[[PHONE_10]]) - hospitalizace pro pneumonii' },\n    @{ Old = '\u2022 Studentka [[PERSON_37]] ([[BIRTH_ID_34]], [[ADDRESS_3]] 602 00) - UK Praha, 2. ro\u010dn\u00edk'; New = '\u2022 Studentka [[PERSON_37]] ([[BIRTH_ID_34]], [[ADDRESS_14]]) - UK Praha, 2. ro\u010dn\u00edk' },\n    @{ Old = '\u2022 Sv\u011bdkyn\u011b PharmDr. [[PERSON_39]] ([[BIRTH_ID_36]], [[PHONE_9]]) - v\u00fdpov\u011b\u010f v trestn\u00ed v\u011bci'; New = '\u2022 Sv\u011bdkyn\u011b PharmDr. [[PERSON_39]] ([[BIRTH_ID_36]], [[PHONE_8]]) - v\u00fdpov\u011b\u010f v trestn\u00ed v\u011bci' },\n)\n\n$appliedCount = 0\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($rep.Old, $false, $true, $false, $false, $false, $true, 1, $false, $rep.New, 2)\n    if ($result) { $appliedCount++ }\n}\n\nWrite-Output \"applied=$appliedCount\""}
